# Push physical sign excel
# The "PhysicalSign" sheet had an extra/duplicate "Height(cm) (cm)" / "162"
# row removed (row 6), and that sheet was left as the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhysicalSign")
$ws.Activate()

# Select + delete row 6 (the duplicate "Height(cm) (cm)" / "162" row)
$ws.Rows("6:6").Select() | Out-Null
$ws.Rows("6:6").Delete() | Out-Null
